# Update 'want to go' counts (column F) per commit 456a3b4 (gh-pages data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 321
$ws.Cells.Item(4, 6).Value = 2994
$ws.Cells.Item(7, 6).Value = 2329
$ws.Cells.Item(8, 6).Value = 1702
$ws.Cells.Item(10, 6).Value = 860
$ws.Cells.Item(11, 6).Value = 128
$ws.Cells.Item(12, 6).Value = 17
$ws.Cells.Item(15, 6).Value = 1544
$ws.Cells.Item(16, 6).Value = 7127
$ws.Cells.Item(18, 6).Value = 7273
$ws.Cells.Item(21, 6).Value = 5542
$ws.Cells.Item(22, 6).Value = 3127
$ws.Cells.Item(23, 6).Value = 3495
$ws.Cells.Item(24, 6).Value = 4
$ws.Cells.Item(26, 6).Value = 192
$ws.Cells.Item(27, 6).Value = 1913
$ws.Cells.Item(29, 6).Value = 307
$ws.Cells.Item(30, 6).Value = 881
$ws.Cells.Item(33, 6).Value = 43
$ws.Cells.Item(34, 6).Value = 2438
$ws.Cells.Item(35, 6).Value = 1229
$ws.Cells.Item(36, 6).Value = 2778
$ws.Cells.Item(37, 6).Value = 38
$ws.Cells.Item(38, 6).Value = 21
$ws.Cells.Item(41, 6).Value = 1097
$ws.Cells.Item(44, 6).Value = 535
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 38
$ws.Cells.Item(8, 6).Value = 219
$ws.Cells.Item(12, 6).Value = 289
$ws.Cells.Item(18, 6).Value = 65
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 76
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 321
$ws.Cells.Item(5, 6).Value = 2994
$ws.Cells.Item(6, 6).Value = 2329
$ws.Cells.Item(7, 6).Value = 1702
$ws.Cells.Item(9, 6).Value = 860
$ws.Cells.Item(10, 6).Value = 128
$ws.Cells.Item(11, 6).Value = 17
$ws.Cells.Item(12, 6).Value = 38
$ws.Cells.Item(13, 6).Value = 76
$ws.Cells.Item(15, 6).Value = 1544
$ws.Cells.Item(16, 6).Value = 219
$ws.Cells.Item(19, 6).Value = 7127
$ws.Cells.Item(21, 6).Value = 7273
$ws.Cells.Item(23, 6).Value = 5542
$ws.Cells.Item(24, 6).Value = 3127
$ws.Cells.Item(25, 6).Value = 289
$ws.Cells.Item(26, 6).Value = 3495
$ws.Cells.Item(30, 6).Value = 1913
$ws.Cells.Item(33, 6).Value = 307
$ws.Cells.Item(34, 6).Value = 881
$ws.Cells.Item(37, 6).Value = 43
$ws.Cells.Item(38, 6).Value = 2438
$ws.Cells.Item(39, 6).Value = 1229
$ws.Cells.Item(40, 6).Value = 65
$ws.Cells.Item(41, 6).Value = 2778
$ws.Cells.Item(42, 6).Value = 38
$ws.Cells.Item(43, 6).Value = 21
$ws.Cells.Item(46, 6).Value = 1097
$ws.Cells.Item(49, 6).Value = 535
